$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by two
# days (45186 -> 45188) for every data row (rows 2 through 302).
$ws.Range("C2:C302").Value = 45188
